$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Formula = "'42.799.17"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.60%  '

# Row 3
$ws.Range("D3").Formula = "'2.322.43"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.76%  '

# Row 4
$ws.Range("E4").Value = '  +0.02%  '

# Row 5
$ws.Range("D5").Formula = "'301.53"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.50%  '

# Row 6
$ws.Range("D6").Formula = "'95.54"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.99%  '

# Row 7
$ws.Range("E7").Value = '  -0.17%  '

# Row 8
$ws.Range("E8").Value = '  +0.11%  '

# Row 9
$ws.Range("D9").Formula = "'0.492"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.86%  '

# Row 10
$ws.Range("D10").Formula = "'34.21"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.75%  '

# Row 11
$ws.Range("D11").Formula = "'18.98"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.37%  '

# Row 12
$ws.Range("E12").Value = '  -0.40%  '

# Row 13
$ws.Range("E13").Value = '  +0.62%  '

# Row 14
$ws.Range("E14").Value = '  -2.25%  '

# Row 15
$ws.Range("D15").Formula = "'2.683.96"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.77%  '

# Row 16
$ws.Range("D16").Formula = "'2.323.20"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.53%  '

# Row 17
$ws.Range("D17").Formula = "'0.789"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.87%  '

# Row 18
$ws.Range("D18").Formula = "'42.756.99"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.42%  '

# Row 19
$ws.Range("D19").Formula = "'12.21"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -4.28%  '

# Row 20
$ws.Range("E20").Value = '  +1.69%  '

# Row 21
$ws.Range("E21").Value = '  -0.84%  '

# Row 22
$ws.Range("D22").Formula = "'67.81"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.68%  '

# Row 23
$ws.Range("D23").Formula = "'2.27"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +5.01%  '

# Row 24
$ws.Range("D24").Formula = "'235.31"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.67%  '

# Row 25
$ws.Range("E25").Value = '  +0.02%  '

# Row 26
$ws.Range("D26").Formula = "'2.40"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.84%  '

# Row 27
$ws.Range("D27").Formula = "'24.41"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.10%  '

# Row 28
$ws.Range("E28").Value = '  +14.37%  '

# Row 29
$ws.Range("E29").Value = '  +0.82%  '

# Row 30
$ws.Range("D30").Formula = "'32.18"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.90%  '

# Row 31
$ws.Range("D31").Formula = "'146.76"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -12.06%  '

# Row 32
$ws.Range("E32").Value = '  -0.04%  '

# Row 33
$ws.Range("E33").Value = '  +0.04%  '

# Row 34
$ws.Range("D34").Formula = "'17.76"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.30%  '

# Row 35
$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").Formula = "'0.0701"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.72%  '

# Row 36
$ws.Range("B36").Value = 'RenderToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D36").Formula = "'4.41"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.12%  '

# Row 37
$ws.Range("E37").Value = '  -1.20%  '

# Row 38
$ws.Range("E38").Value = '  +3.14%  '

# Row 39
$ws.Range("E39").Value = '  -1.06%  '

# Row 40
$ws.Range("D40").Formula = "'2.75"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.53%  '

# Row 41
$ws.Range("D41").Formula = "'22.09"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +22.93%  '

# Row 42
$ws.Range("E42").Value = '  -0.96%  '

# Row 43
$ws.Range("D43").Formula = "'1.924.92"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.64%  '

# Row 44
$ws.Range("E44").Value = '  -1.07%  '

# Row 45
$ws.Range("D45").Formula = "'10.12"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.24%  '

# Row 46
$ws.Range("E46").Value = '  -1.76%  '

# Row 47
$ws.Range("D47").Formula = "'2.74"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.54%  '

# Row 48
$ws.Range("E48").Value = '  -0.74%  '

# Row 49
$ws.Range("D49").Formula = "'2.552.04"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.82%  '

# Row 50
$ws.Range("E50").Value = '  -0.39%  '

# Row 51
$ws.Range("D51").Formula = "'72.26"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.05%  '
